$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 442, shifting existing rows 442:539 down to 443:540
$ws.Rows("442:442").Insert()

# Fill the new row 442 with values:
# category/classification columns copied from the row now at 443 (the old row 442 data)
$ws.Range("A442").Value = $ws.Range("A443").Value()
$ws.Range("B442").Value = $ws.Range("B443").Value()
$ws.Range("C442").Value = $ws.Range("C443").Value()
$ws.Range("D442").Value = 45173
$ws.Range("E442").Value = $ws.Range("E443").Value()
$ws.Range("F442").Value = $ws.Range("F443").Value()
$ws.Range("G442").Value = $ws.Range("G443").Value()
$ws.Range("H442").Value = $ws.Range("H443").Value()
$ws.Range("I442").Value = $ws.Range("I443").Value()
$ws.Range("J442").Value = 300
$ws.Range("K442").Value = 20000
$ws.Range("L442").Value = 20000
$ws.Range("M442").Value = 20000
$ws.Range("N442").Value = $ws.Range("N443").Value()
$ws.Range("O442").Value = $ws.Range("O443").Value()
$ws.Range("P442").Value = 2000
$ws.Range("Q442").Value = $ws.Range("Q443").Value()
$ws.Range("R442").Value = $ws.Range("R443").Value()
